$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows (continuing the series through 13/05/2021, serials 44326-44329)
$data = @(
    @(44326, 1, 2, 32.34675723758694),
    @(44327, 0, 2, 32.34675723758694),
    @(44328, 1, 3, 48.5201358563804),
    @(44329, 0, 3, 48.5201358563804)
)

$lastRow = 251
$r = $lastRow + 1
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]

    # Match the date-styled format used by the rest of column A
    $ws.Range("A" + $lastRow).Copy()
    $ws.Range("A" + $r).PasteSpecial(-4122)

    $r = $r + 1
}

$excel.CutCopyMode = 0
